# 7.1.2 workbook update: add 2021 data point, bump copyright year to 2022.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")

# --- 1. Tabelle2 (chart source data): append the 2021 row -------------------
$ws2.Range("A12").Value = 2021
$ws2.Range("B12").Value = 95
$ws2.Range("C12").Value = 5

# --- 2. Tabelle1: make room for the new year row by inserting a row --------
# (row 16 is already a blank placeholder row that receives the 2021 data;
#  a new blank row is inserted at 17 so the footer block below shifts down
#  by one, matching how the sheet grows every year.)
$ws1.Rows.Item(17).Insert()

$ws1.Range("A16").Value = 2021
$ws1.Range("B16").Value = "> 95"

# --- 3. Update the copyright year in the footer (now row 39) ---------------
$ws1.Range("B39").Value = "©       Statistisches Bundesamt (Destatis) 2022"

# --- 4. Reposition the chart + small picture that sit below the data -------
# The chart (graphicFrame) anchor and the small picture anchor are both
# pinned to absolute rows, so shifting the footer down means they need to
# move down by one row as well.
$co = $ws1.ChartObjects().Item(1)
$targetTop    = $ws1.Rows.Item(18).Top
$targetLeft   = $ws1.Columns.Item(1).Left
$targetBottom = $ws1.Rows.Item(36).Top + (119743 / 12700)
$targetRight  = $ws1.Columns.Item(8).Left + (283028 / 12700)
$co.Top    = $targetTop
$co.Left   = $targetLeft
$co.Width  = $targetRight - $targetLeft
$co.Height = $targetBottom - $targetTop

$pic = $ws1.Shapes.Item("Grafik 6")
$pic.Top = $ws1.Rows.Item(39).Top + (41412 / 12700)

# --- 5. Keep the chart series bound to the full (now 12-row) range ---------
$chart = $co.Chart
$s1 = $chart.SeriesCollection().Item(1)
$s1.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$12,Tabelle2!`$B`$1:`$B`$12,1)"
$s2 = $chart.SeriesCollection().Item(2)
$s2.Formula = "=SERIES(,Tabelle2!`$A`$1:`$A`$12,Tabelle2!`$C`$1:`$C`$12,2)"
